# Apply updated crypto price/volume data (and two row swaps) per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.587.64'
$ws.Range("E2").Value = '  +4.59%  '

$ws.Range("D3").Value = '2.266.96'
$ws.Range("E3").Value = '  +3.65%  '

$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("E5").Value = '  +3.47%  '

$ws.Range("D6").Value = '''92.76'
$ws.Range("E6").Value = '  +7.18%  '

$ws.Range("D7").Value = '''0.523'
$ws.Range("E7").Value = '  +3.13%  '

$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("D9").Value = '''0.488'
$ws.Range("E9").Value = '  +5.00%  '

$ws.Range("D10").Value = '''54.11'
$ws.Range("E10").Value = '  +8.40%  '

$ws.Range("D11").Value = '''32.26'
$ws.Range("E11").Value = '  +8.22%  '

$ws.Range("D12").Value = '''0.0797'
$ws.Range("E12").Value = '  +3.10%  '

$ws.Range("E13").Value = '  +3.51%  '

$ws.Range("D14").Value = '''6.65'
$ws.Range("E14").Value = '  +4.01%  '

$ws.Range("D15").Value = '2.618.64'
$ws.Range("E15").Value = '  +3.63%  '

$ws.Range("D16").Value = '''14.24'
$ws.Range("E16").Value = '  +4.84%  '

$ws.Range("D17").Value = '2.285.59'
$ws.Range("E17").Value = '  +5.89%  '

$ws.Range("D18").Value = '''0.756'
$ws.Range("E18").Value = '  +5.07%  '

$ws.Range("D19").Value = '41.496.27'
$ws.Range("E19").Value = '  +4.67%  '

$ws.Range("D20").Value = '''12.46'
$ws.Range("E20").Value = '  +12.46%  '

$ws.Range("D21").Value = '0.0₃0909'
$ws.Range("E21").Value = '  +3.54%  '

$ws.Range("D22").Value = '''5.94'
$ws.Range("E22").Value = '  +4.05%  '

$ws.Range("D23").Value = '''67.11'
$ws.Range("E23").Value = '  +3.37%  '

$ws.Range("D24").Value = '''240.66'
$ws.Range("E24").Value = '  +2.40%  '

$ws.Range("D25").Value = '''2.58'
$ws.Range("E25").Value = '  +5.54%  '

$ws.Range("E26").Value = '  -0.16%  '

$ws.Range("E27").Value = '  +4.84%  '

$ws.Range("D28").Value = '''23.82'
$ws.Range("E28").Value = '  +6.51%  '

$ws.Range("E29").Value = '  +1.99%  '

$ws.Range("D30").Value = '''9.71'
$ws.Range("E30").Value = '  +6.89%  '

$ws.Range("D31").Value = '''34.07'
$ws.Range("E31").Value = '  +9.74%  '

$ws.Range("D32").Value = '''157.39'
$ws.Range("E32").Value = '  +1.45%  '

$ws.Range("D33").Value = '''0.999'
$ws.Range("E33").Value = '  +0.05%  '

$ws.Range("D34").Value = '''5.21'
$ws.Range("E34").Value = '  +7.51%  '

$ws.Range("D35").Value = '''0.0738'
$ws.Range("E35").Value = '  +4.96%  '

$ws.Range("E36").Value = '  +9.71%  '

$ws.Range("E37").Value = '  +1.58%  '

$ws.Range("D38").Value = '''16.62'
$ws.Range("E38").Value = '  +10.11%  '

$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").Value = '''0.104'
$ws.Range("E39").Value = '  +7.63%  '

$ws.Range("B40").Value = 'Stellar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D40").Value = '''0.115'
$ws.Range("E40").Value = '  +3.13%  '

$ws.Range("D41").Value = '''1.78'
$ws.Range("E41").Value = '  +7.70%  '

$ws.Range("D42").Value = '''4.02'
$ws.Range("E42").Value = '  +8.43%  '

$ws.Range("D43").Value = '''20.62'
$ws.Range("E43").Value = '  +20.86%  '

$ws.Range("D44").Value = '2.067.46'
$ws.Range("E44").Value = '  -2.05%  '

$ws.Range("D45").Value = '''0.0278'
$ws.Range("E45").Value = '  +4.78%  '

$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D46").Value = '''10.11'
$ws.Range("E46").Value = '  +4.94%  '

$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").Value = '''2.98'
$ws.Range("E47").Value = '  +13.18%  '

$ws.Range("D48").Value = '''2.00'
$ws.Range("E48").Value = '  -5.10%  '

$ws.Range("D49").Value = '2.486.61'
$ws.Range("E49").Value = '  +3.69%  '

$ws.Range("D50").Value = '''1.53'
$ws.Range("E50").Value = '  +4.57%  '

$ws.Range("D51").Value = '''1.15'
$ws.Range("E51").Value = '  +5.27%  '
